# Update "hasil tabulasi" worksheet with refreshed tabulation numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Real Number Jenis Kelamin (A1 block)
$ws.Range("C3").Value = 234
$ws.Range("C4").Value = 266

# Percentage Jenis Kelamin (A7 block)
$ws.Range("C9").Value = 46.8
$ws.Range("C10").Value = 53.2

# Real Number and Percentage Jenis Kelamin (A13 block)
$ws.Range("C15").Value = 234
$ws.Range("C16").Value = 266
$ws.Range("C18").Value = 46.8
$ws.Range("C19").Value = 53.2

# gender vs usia counts (row 22 block)
$ws.Range("C23").Value = 234
$ws.Range("C24").Value = 266
$ws.Range("C26").Value = 19
$ws.Range("C27").Value = 35
$ws.Range("C28").Value = 84
$ws.Range("C29").Value = 37
$ws.Range("C30").Value = 325

# Crosstabulasi gender vs usia (A33 block)
$ws.Range("C36").Value = 5.55555555555556
$ws.Range("D36").Value = 2.25563909774436
$ws.Range("E36").Value = 3.8
$ws.Range("C37").Value = 6.83760683760684
$ws.Range("D37").Value = 7.14285714285714
$ws.Range("E37").Value = 7
$ws.Range("C38").Value = 17.9487179487179
$ws.Range("D38").Value = 15.7894736842105
$ws.Range("E38").Value = 16.8
$ws.Range("C39").Value = 7.69230769230769
$ws.Range("D39").Value = 7.14285714285714
$ws.Range("E39").Value = 7.4
$ws.Range("C40").Value = 61.965811965812
$ws.Range("D40").Value = 67.6691729323308
$ws.Range("E40").Value = 65
$ws.Range("C41").Value = 234
$ws.Range("D41").Value = 266

# Crosstabulasi gender vs awareness (A43 block)
$ws.Range("C46").Value = 48.546511627907
$ws.Range("D46").Value = 42.9487179487179
$ws.Range("E46").Value = 46.8
$ws.Range("C47").Value = 51.453488372093
$ws.Range("D47").Value = 57.0512820512821
$ws.Range("E47").Value = 53.2
$ws.Range("C48").Value = 344
$ws.Range("D48").Value = 156

# Percentage Stasiun TV (A50 block)
$ws.Range("B52").Value = 91.1646586345382
$ws.Range("B53").Value = 69.4779116465863
$ws.Range("B54").Value = 49.1967871485944
$ws.Range("B55").Value = 28.1124497991968
$ws.Range("B56").Value = 1.40562248995984
$ws.Range("B57").Value = 498

# Crosstabulasi stasiun TV vs gender (A59 block)
$ws.Range("B62").Value = 92.7038626609442
$ws.Range("C62").Value = 89.811320754717
$ws.Range("D62").Value = 91.1646586345382
$ws.Range("B63").Value = 66.9527896995708
$ws.Range("C63").Value = 71.6981132075472
$ws.Range("D63").Value = 69.4779116465863
$ws.Range("B64").Value = 47.2103004291845
$ws.Range("C64").Value = 50.9433962264151
$ws.Range("D64").Value = 49.1967871485944
$ws.Range("B65").Value = 31.7596566523605
$ws.Range("C65").Value = 24.9056603773585
$ws.Range("D65").Value = 28.1124497991968
$ws.Range("B66").Value = 0.858369098712446
$ws.Range("C66").Value = 1.88679245283019
$ws.Range("D66").Value = 1.40562248995984
$ws.Range("B67").Value = 233
$ws.Range("C67").Value = 265
$ws.Range("D67").Value = 498

# Crosstabulasi stasiun TV vs gender, gender vs usia (A69 block)
$ws.Range("C73").Value = 62
$ws.Range("D73").Value = 67.7
$ws.Range("E73").Value = 65
$ws.Range("C74").Value = 17.9
$ws.Range("D74").Value = 15.8
$ws.Range("E74").Value = 16.8

# Rows 75 and 76 swap their "usia" category label along with refreshed values
$ws.Range("B75").Value = "26 - 30 th"
$ws.Range("C75").Value = 7.7
$ws.Range("D75").Value = 7.1
$ws.Range("E75").Value = 7.4

$ws.Range("B76").Value = "16 - 20 th"
$ws.Range("C76").Value = 6.8
$ws.Range("D76").Value = 7.1
$ws.Range("E76").Value = 7

$ws.Range("C77").Value = 5.6
$ws.Range("D77").Value = 2.3
$ws.Range("E77").Value = 3.8
$ws.Range("C78").Value = 234
$ws.Range("D78").Value = 266
